$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sdfewf")
$ws.Name = "nadibf"

$ws.Cells.Item(1, 10).Value = 42.04670763015747
$ws.Cells.Item(2, 10).Value = 40.82866406440735
$ws.Cells.Item(3, 10).Value = 40.51431345939636
$ws.Cells.Item(4, 10).Value = 42.34519624710083
$ws.Cells.Item(5, 2).Value = 2235
$ws.Cells.Item(5, 4).Value = 2224
$ws.Cells.Item(5, 5).Value = 10
$ws.Cells.Item(5, 6).Value = 3
$ws.Cells.Item(5, 7).Value = 99.86528962730131
$ws.Cells.Item(5, 8).Value = 99.55237242614145
$ws.Cells.Item(5, 9).Value = 0.005834829443447037
$ws.Cells.Item(5, 10).Value = 41.70451474189758
$ws.Cells.Item(6, 2).Value = 2584
$ws.Cells.Item(6, 5).Value = 27
$ws.Cells.Item(6, 8).Value = 98.95470383275261
$ws.Cells.Item(6, 9).Value = 0.01594710229482691
$ws.Cells.Item(6, 10).Value = 41.18650507926941
$ws.Cells.Item(7, 2).Value = 2026
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 8).Value = 99.90123456790124
$ws.Cells.Item(7, 9).Value = 0.00246669955599408
$ws.Cells.Item(7, 10).Value = 43.21181845664978
$ws.Cells.Item(8, 2).Value = 2222
$ws.Cells.Item(8, 5).Value = 88
$ws.Cells.Item(8, 8).Value = 96.03782080144079
$ws.Cells.Item(8, 9).Value = 0.04168618266978923
$ws.Cells.Item(8, 10).Value = 42.3247447013855
$ws.Cells.Item(9, 2).Value = 1768
$ws.Cells.Item(9, 5).Value = 13
$ws.Cells.Item(9, 8).Value = 99.26428975664969
$ws.Cells.Item(9, 9).Value = 0.01135073779795687
$ws.Cells.Item(9, 10).Value = 39.94191884994507
$ws.Cells.Item(10, 10).Value = 39.85851001739502
$ws.Cells.Item(11, 10).Value = 40.13940715789795
$ws.Cells.Item(12, 2).Value = 2559
$ws.Cells.Item(12, 5).Value = 21
$ws.Cells.Item(12, 8).Value = 99.1790461297889
$ws.Cells.Item(12, 9).Value = 0.008274231678486997
$ws.Cells.Item(12, 10).Value = 39.72089266777039
$ws.Cells.Item(13, 2).Value = 1812
$ws.Cells.Item(13, 5).Value = 18
$ws.Cells.Item(13, 8).Value = 99.00607399226946
$ws.Cells.Item(13, 9).Value = 0.0105849582172702
$ws.Cells.Item(13, 10).Value = 39.43358469009399
$ws.Cells.Item(14, 10).Value = 39.27301263809204
$ws.Cells.Item(15, 10).Value = 38.85359668731689
$ws.Cells.Item(16, 10).Value = 40.73881936073303
$ws.Cells.Item(17, 10).Value = 39.66579484939575
$ws.Cells.Item(18, 10).Value = 41.65223956108093
$ws.Cells.Item(19, 2).Value = 1992
$ws.Cells.Item(19, 5).Value = 5
$ws.Cells.Item(19, 8).Value = 99.74886991461577
$ws.Cells.Item(19, 9).Value = 0.002516356316054353
$ws.Cells.Item(19, 10).Value = 40.30982542037964
$ws.Cells.Item(20, 10).Value = 39.63760089874268
$ws.Cells.Item(21, 10).Value = 40.67124652862549
$ws.Cells.Item(22, 2).Value = 1520
$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 8).Value = 99.80250164581962
$ws.Cells.Item(22, 9).Value = 0.001977587343441002
$ws.Cells.Item(22, 10).Value = 39.27433133125305
$ws.Cells.Item(23, 10).Value = 39.93972229957581
$ws.Cells.Item(24, 2).Value = 2601
$ws.Cells.Item(24, 4).Value = 2598
$ws.Cells.Item(24, 6).Value = 2
$ws.Cells.Item(24, 7).Value = 99.92307692307692
$ws.Cells.Item(24, 8).Value = 99.92307692307692
$ws.Cells.Item(24, 9).Value = 0.001537870049980777
$ws.Cells.Item(24, 10).Value = 40.18035435676575
$ws.Cells.Item(25, 2).Value = 1951
$ws.Cells.Item(25, 4).Value = 1949
$ws.Cells.Item(25, 6).Value = 13
$ws.Cells.Item(25, 7).Value = 99.33741080530072
$ws.Cells.Item(25, 8).Value = 99.94871794871794
$ws.Cells.Item(25, 9).Value = 0.007131940906775344
$ws.Cells.Item(25, 10).Value = 39.88451290130615
$ws.Cells.Item(26, 10).Value = 41.18431997299194
$ws.Cells.Item(27, 2).Value = 2883
$ws.Cells.Item(27, 4).Value = 2868
$ws.Cells.Item(27, 5).Value = 14
$ws.Cells.Item(27, 6).Value = 110
$ws.Cells.Item(27, 7).Value = 96.30624580255204
$ws.Cells.Item(27, 8).Value = 99.51422623178348
$ws.Cells.Item(27, 9).Value = 0.04162470627727426
$ws.Cells.Item(27, 10).Value = 40.62373661994934
$ws.Cells.Item(28, 2).Value = 2636
$ws.Cells.Item(28, 4).Value = 2635
$ws.Cells.Item(28, 6).Value = 20
$ws.Cells.Item(28, 7).Value = 99.24670433145009
$ws.Cells.Item(28, 9).Value = 0.007530120481927711
$ws.Cells.Item(28, 10).Value = 39.78992247581482
$ws.Cells.Item(29, 2).Value = 2086
$ws.Cells.Item(29, 3).Value = 2331
$ws.Cells.Item(29, 4).Value = 2075
$ws.Cells.Item(29, 5).Value = 10
$ws.Cells.Item(29, 6).Value = 255
$ws.Cells.Item(29, 7).Value = 89.05579399141631
$ws.Cells.Item(29, 8).Value = 99.52038369304556
$ws.Cells.Item(29, 9).Value = 0.1136851136851137
$ws.Cells.Item(29, 10).Value = 39.48236966133118
$ws.Cells.Item(30, 2).Value = 2945
$ws.Cells.Item(30, 4).Value = 2936
$ws.Cells.Item(30, 5).Value = 8
$ws.Cells.Item(30, 6).Value = 16
$ws.Cells.Item(30, 7).Value = 99.4579945799458
$ws.Cells.Item(30, 8).Value = 99.72826086956522
$ws.Cells.Item(30, 9).Value = 0.008127328140873687
$ws.Cells.Item(30, 10).Value = 40.15036797523499
$ws.Cells.Item(31, 2).Value = 3005
$ws.Cells.Item(31, 4).Value = 3004
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 100
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 40.37878513336182
$ws.Cells.Item(32, 2).Value = 2625
$ws.Cells.Item(32, 4).Value = 2622
$ws.Cells.Item(32, 5).Value = 2
$ws.Cells.Item(32, 6).Value = 26
$ws.Cells.Item(32, 7).Value = 99.01812688821752
$ws.Cells.Item(32, 8).Value = 99.92378048780488
$ws.Cells.Item(32, 9).Value = 0.01057002642506606
$ws.Cells.Item(32, 10).Value = 39.32681894302368
$ws.Cells.Item(33, 10).Value = 44.94929528236389
$ws.Cells.Item(34, 10).Value = 45.38955593109131
$ws.Cells.Item(35, 2).Value = 2257
$ws.Cells.Item(35, 4).Value = 2256
$ws.Cells.Item(35, 6).Value = 4
$ws.Cells.Item(35, 7).Value = 99.82300884955752
$ws.Cells.Item(35, 9).Value = 0.001769128704113224
$ws.Cells.Item(35, 10).Value = 42.1408371925354
$ws.Cells.Item(36, 2).Value = 3358
$ws.Cells.Item(36, 4).Value = 3357
$ws.Cells.Item(36, 6).Value = 4
$ws.Cells.Item(36, 7).Value = 99.88098780124963
$ws.Cells.Item(36, 9).Value = 0.001189767995240928
$ws.Cells.Item(36, 10).Value = 43.1296055316925
$ws.Cells.Item(37, 2).Value = 2253
$ws.Cells.Item(37, 5).Value = 49
$ws.Cells.Item(37, 8).Value = 97.82415630550622
$ws.Cells.Item(37, 9).Value = 0.0240036231884058
$ws.Cells.Item(37, 10).Value = 42.35402250289917
$ws.Cells.Item(38, 10).Value = 40.29392242431641
$ws.Cells.Item(39, 10).Value = 40.13136625289917
$ws.Cells.Item(40, 10).Value = 40.02736830711365
$ws.Cells.Item(41, 2).Value = 2486
$ws.Cells.Item(41, 4).Value = 2480
$ws.Cells.Item(41, 6).Value = 1
$ws.Cells.Item(41, 7).Value = 99.95969367190649
$ws.Cells.Item(41, 8).Value = 99.79879275653923
$ws.Cells.Item(41, 9).Value = 0.0024174053182917
$ws.Cells.Item(41, 10).Value = 40.13996267318726
$ws.Cells.Item(42, 10).Value = 40.51850080490112
$ws.Cells.Item(43, 2).Value = 2054
$ws.Cells.Item(43, 5).Value = 9
$ws.Cells.Item(43, 8).Value = 99.56161714564053
$ws.Cells.Item(43, 9).Value = 0.007797270955165692
$ws.Cells.Item(43, 10).Value = 40.27977967262268
$ws.Cells.Item(44, 2).Value = 2256
$ws.Cells.Item(44, 5).Value = 1
$ws.Cells.Item(44, 8).Value = 99.95565410199556
$ws.Cells.Item(44, 9).Value = 0.0004434589800443459
$ws.Cells.Item(44, 10).Value = 40.15129661560059
$ws.Cells.Item(45, 2).Value = 1572
$ws.Cells.Item(45, 5).Value = 2
$ws.Cells.Item(45, 8).Value = 99.87269255251432
$ws.Cells.Item(45, 9).Value = 0.001273885350318471
$ws.Cells.Item(45, 10).Value = 39.85698246955872
$ws.Cells.Item(46, 2).Value = 1786
$ws.Cells.Item(46, 5).Value = 6
$ws.Cells.Item(46, 8).Value = 99.66386554621849
$ws.Cells.Item(46, 9).Value = 0.003370786516853933
$ws.Cells.Item(46, 10).Value = 40.05977582931519
$ws.Cells.Item(47, 2).Value = 3076
$ws.Cells.Item(47, 4).Value = 3074
$ws.Cells.Item(47, 6).Value = 2
$ws.Cells.Item(47, 7).Value = 99.93498049414825
$ws.Cells.Item(47, 8).Value = 99.96747967479675
$ws.Cells.Item(47, 9).Value = 0.0009749756256093598
$ws.Cells.Item(47, 10).Value = 42.01871585845947
$ws.Cells.Item(48, 10).Value = 41.22881722450256
